# Scheduled runner update: refresh Leve profit calculations (currentAveragePrice,
# LevePriceNQ/HQ, LeveProfitNQ/HQ, etc.) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(17, 8).Value = 4173205.5  # H17
$ws.Cells.Item(17, 10).Value = 4354567  # J17
$ws.Cells.Item(17, 12).Value = 13063701  # L17
$ws.Cells.Item(17, 14).Value = -13064037  # N17
$ws.Cells.Item(51, 8).Value = 7000.3335  # H51
$ws.Cells.Item(51, 10).Value = 4500.5  # J51
$ws.Cells.Item(51, 12).Value = 4500.5  # L51
$ws.Cells.Item(51, 14).Value = -5468.5  # N51
$ws.Cells.Item(88, 8).Value = 1345  # H88
$ws.Cells.Item(88, 9).Value = 503  # I88
$ws.Cells.Item(88, 10).Value = 1450.25  # J88
$ws.Cells.Item(88, 11).Value = 503  # K88
$ws.Cells.Item(88, 12).Value = 1450.25  # L88
$ws.Cells.Item(88, 13).Value = -97  # M88
$ws.Cells.Item(88, 14).Value = -2262.25  # N88
$ws.Cells.Item(91, 8).Value = 1345  # H91
$ws.Cells.Item(91, 9).Value = 503  # I91
$ws.Cells.Item(91, 10).Value = 1450.25  # J91
$ws.Cells.Item(91, 11).Value = 503  # K91
$ws.Cells.Item(91, 12).Value = 1450.25  # L91
$ws.Cells.Item(91, 13).Value = 901  # M91
$ws.Cells.Item(91, 14).Value = -4258.25  # N91
$ws.Cells.Item(129, 8).Value = 228157.64  # H129
$ws.Cells.Item(129, 10).Value = 286759.12  # J129
$ws.Cells.Item(129, 12).Value = 860277.36  # L129
$ws.Cells.Item(129, 14).Value = -870277.36  # N129
$ws.Cells.Item(132, 8).Value = 20410186  # H132
$ws.Cells.Item(132, 10).Value = 0  # J132
$ws.Cells.Item(132, 12).Value = 0  # L132
$ws.Cells.Item(132, 14).ClearContents()  # N132
$ws.Cells.Item(138, 8).Value = 3571.23  # H138
$ws.Cells.Item(138, 9).Value = 2424.125  # I138
$ws.Cells.Item(138, 10).Value = 3789.726  # J138
$ws.Cells.Item(138, 11).Value = 7272.375  # K138
$ws.Cells.Item(138, 12).Value = 11369.178  # L138
$ws.Cells.Item(138, 13).Value = -2132.375  # M138
$ws.Cells.Item(138, 14).Value = -21649.178  # N138
$ws.Cells.Item(139, 8).Value = 44577.145  # H139
$ws.Cells.Item(139, 10).Value = 50680  # J139
$ws.Cells.Item(139, 12).Value = 50680  # L139
$ws.Cells.Item(139, 14).Value = -60960  # N139
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(45, 8).Value = 5425.3076  # H45
$ws.Cells.Item(45, 9).Value = 7080.2  # I45
$ws.Cells.Item(45, 10).Value = 4391  # J45
$ws.Cells.Item(45, 11).Value = 7080.2  # K45
$ws.Cells.Item(45, 12).Value = 4391  # L45
$ws.Cells.Item(45, 13).Value = -6703.2  # M45
$ws.Cells.Item(45, 14).Value = -5145  # N45
$ws.Cells.Item(133, 8).Value = 69666.664  # H133
$ws.Cells.Item(133, 9).Value = 14000  # I133
$ws.Cells.Item(133, 11).Value = 14000  # K133
$ws.Cells.Item(133, 13).Value = -11470  # M133
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(20, 8).Value = 3564.2856  # H20
$ws.Cells.Item(20, 9).Value = 3376.923  # I20
$ws.Cells.Item(20, 10).Value = 6000  # J20
$ws.Cells.Item(20, 11).Value = 3376.923  # K20
$ws.Cells.Item(20, 12).Value = 6000  # L20
$ws.Cells.Item(20, 13).Value = -3129.923  # M20
$ws.Cells.Item(20, 14).Value = -6494  # N20
$ws.Cells.Item(86, 8).Value = 1918.3214  # H86
$ws.Cells.Item(86, 9).Value = 1790.8096  # I86
$ws.Cells.Item(86, 10).Value = 2300.8572  # J86
$ws.Cells.Item(86, 11).Value = 1790.8096  # K86
$ws.Cells.Item(86, 12).Value = 2300.8572  # L86
$ws.Cells.Item(86, 13).Value = -667.8096  # M86
$ws.Cells.Item(86, 14).Value = -4546.8572  # N86
$ws.Cells.Item(89, 8).Value = 1918.3214  # H89
$ws.Cells.Item(89, 9).Value = 1790.8096  # I89
$ws.Cells.Item(89, 10).Value = 2300.8572  # J89
$ws.Cells.Item(89, 11).Value = 8954.048000000001  # K89
$ws.Cells.Item(89, 12).Value = 11504.286  # L89
$ws.Cells.Item(89, 13).Value = -3338.048000000001  # M89
$ws.Cells.Item(89, 14).Value = -22736.286  # N89
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(31, 8).Value = 6936.25  # H31
$ws.Cells.Item(31, 9).Value = 3495  # I31
$ws.Cells.Item(31, 10).Value = 8353.235000000001  # J31
$ws.Cells.Item(31, 11).Value = 3495  # K31
$ws.Cells.Item(31, 12).Value = 8353.235000000001  # L31
$ws.Cells.Item(31, 13).Value = -3200  # M31
$ws.Cells.Item(31, 14).Value = -8943.235000000001  # N31
$ws.Cells.Item(34, 8).Value = 6936.25  # H34
$ws.Cells.Item(34, 9).Value = 3495  # I34
$ws.Cells.Item(34, 10).Value = 8353.235000000001  # J34
$ws.Cells.Item(34, 11).Value = 3495  # K34
$ws.Cells.Item(34, 12).Value = 8353.235000000001  # L34
$ws.Cells.Item(34, 13).Value = -3293  # M34
$ws.Cells.Item(34, 14).Value = -8757.235000000001  # N34
$ws.Cells.Item(99, 8).Value = 3672.9285  # H99
$ws.Cells.Item(99, 9).Value = 2716.2856  # I99
$ws.Cells.Item(99, 10).Value = 6542.857  # J99
$ws.Cells.Item(99, 11).Value = 2716.2856  # K99
$ws.Cells.Item(99, 12).Value = 6542.857  # L99
$ws.Cells.Item(99, 13).Value = -1218.2856  # M99
$ws.Cells.Item(99, 14).Value = -9538.857  # N99
$ws.Cells.Item(122, 8).Value = 4585.125  # H122
$ws.Cells.Item(122, 10).Value = 2795  # J122
$ws.Cells.Item(122, 12).Value = 8385  # L122
$ws.Cells.Item(122, 14).Value = -13285  # N122
$ws.Cells.Item(126, 8).Value = 3672.9285  # H126
$ws.Cells.Item(126, 9).Value = 2716.2856  # I126
$ws.Cells.Item(126, 10).Value = 6542.857  # J126
$ws.Cells.Item(126, 11).Value = 8148.8568  # K126
$ws.Cells.Item(126, 12).Value = 19628.571  # L126
$ws.Cells.Item(126, 13).Value = -5678.8568  # M126
$ws.Cells.Item(126, 14).Value = -24568.571  # N126
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(131, 8).Value = 745.64  # H131
$ws.Cells.Item(131, 9).Value = 570  # I131
$ws.Cells.Item(131, 10).Value = 749.2245  # J131
$ws.Cells.Item(131, 11).Value = 1710  # K131
$ws.Cells.Item(131, 12).Value = 2247.6735  # L131
$ws.Cells.Item(131, 13).Value = 3330  # M131
$ws.Cells.Item(131, 14).Value = -12327.6735  # N131
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(59, 8).Value = 18000  # H59
$ws.Cells.Item(59, 10).Value = 18000  # J59
$ws.Cells.Item(59, 12).Value = 18000  # L59
$ws.Cells.Item(59, 14).Value = -19166  # N59
$ws.Cells.Item(70, 8).Value = 2407939.5  # H70
$ws.Cells.Item(70, 9).Value = 4120.5625  # I70
$ws.Cells.Item(70, 11).Value = 4120.5625  # K70
$ws.Cells.Item(70, 13).Value = -3850.5625  # M70
$ws.Cells.Item(73, 8).Value = 2407939.5  # H73
$ws.Cells.Item(73, 9).Value = 4120.5625  # I73
$ws.Cells.Item(73, 11).Value = 4120.5625  # K73
$ws.Cells.Item(73, 13).Value = -3184.5625  # M73
$ws.Cells.Item(80, 8).Value = 3624.08  # H80
$ws.Cells.Item(80, 9).Value = 3263.5454  # I80
$ws.Cells.Item(80, 10).Value = 3907.3572  # J80
$ws.Cells.Item(80, 11).Value = 3263.5454  # K80
$ws.Cells.Item(80, 12).Value = 3907.3572  # L80
$ws.Cells.Item(80, 13).Value = -2265.5454  # M80
$ws.Cells.Item(80, 14).Value = -5903.3572  # N80
$ws.Cells.Item(83, 8).Value = 3624.08  # H83
$ws.Cells.Item(83, 9).Value = 3263.5454  # I83
$ws.Cells.Item(83, 10).Value = 3907.3572  # J83
$ws.Cells.Item(83, 11).Value = 16317.727  # K83
$ws.Cells.Item(83, 12).Value = 19536.786  # L83
$ws.Cells.Item(83, 13).Value = -11325.727  # M83
$ws.Cells.Item(83, 14).Value = -29520.786  # N83
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(16, 8).Value = 548.6875  # H16
$ws.Cells.Item(16, 9).Value = 578.5  # I16
$ws.Cells.Item(16, 11).Value = 578.5  # K16
$ws.Cells.Item(16, 13).Value = -408.5  # M16
$ws.Cells.Item(22, 8).Value = 3175.5  # H22
$ws.Cells.Item(22, 9).Value = 3567.3333  # I22
$ws.Cells.Item(22, 10).Value = 2000  # J22
$ws.Cells.Item(22, 11).Value = 3567.3333  # K22
$ws.Cells.Item(22, 12).Value = 2000  # L22
$ws.Cells.Item(22, 13).Value = -3272.3333  # M22
$ws.Cells.Item(22, 14).Value = -2590  # N22
$ws.Cells.Item(27, 8).Value = 3175.5  # H27
$ws.Cells.Item(27, 9).Value = 3567.3333  # I27
$ws.Cells.Item(27, 10).Value = 2000  # J27
$ws.Cells.Item(27, 11).Value = 3567.3333  # K27
$ws.Cells.Item(27, 12).Value = 2000  # L27
$ws.Cells.Item(27, 13).Value = -3460.3333  # M27
$ws.Cells.Item(27, 14).Value = -2214  # N27
$ws.Cells.Item(68, 8).Value = 2369.8235  # H68
$ws.Cells.Item(68, 9).Value = 1685.5714  # I68
$ws.Cells.Item(68, 11).Value = 1685.5714  # K68
$ws.Cells.Item(68, 13).Value = -936.5714  # M68
$ws.Cells.Item(71, 8).Value = 2369.8235  # H71
$ws.Cells.Item(71, 9).Value = 1685.5714  # I71
$ws.Cells.Item(71, 11).Value = 8427.857  # K71
$ws.Cells.Item(71, 13).Value = -4683.857  # M71
$ws.Cells.Item(122, 8).Value = 1964391.4  # H122
$ws.Cells.Item(122, 9).Value = 2181935  # I122
$ws.Cells.Item(122, 10).Value = 6500  # J122
$ws.Cells.Item(122, 11).Value = 6545805  # K122
$ws.Cells.Item(122, 12).Value = 19500  # L122
$ws.Cells.Item(122, 13).Value = -6543355  # M122
$ws.Cells.Item(122, 14).Value = -24400  # N122
$ws.Cells.Item(132, 8).Value = 3465.1  # H132
$ws.Cells.Item(132, 9).Value = 2243.4285  # I132
$ws.Cells.Item(132, 10).Value = 6315.6665  # J132
$ws.Cells.Item(132, 11).Value = 6730.2855  # K132
$ws.Cells.Item(132, 12).Value = 18946.9995  # L132
$ws.Cells.Item(132, 13).Value = -4200.2855  # M132
$ws.Cells.Item(132, 14).Value = -24006.9995  # N132
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(140, 8).Value = 46825  # H140
$ws.Cells.Item(140, 10).Value = 46825  # J140
$ws.Cells.Item(140, 12).Value = 46825  # L140
$ws.Cells.Item(140, 14).Value = -57185  # N140
